$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new cell C14 (Status = "ID Used") to existing row 14
$ws.Range("C14").Value = "ID Used"

# Add new row 15: IDNumber + Status
$ws.Range("A15").Value = "0105052295182"
$ws.Range("C15").Value = "ID Used"

# Add new row 16: IDNumber only
$ws.Range("A16").Value = "9805051386186"
